# Update "想去人数" (want-to-go count, column F) values on the "展览" and
# "全部类型" sheets to the freshly scraped figures from the gh-pages data
# refresh (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# row -> new F value, for the "展览" worksheet
$sheet1Updates = @{
    2  = 15172
    3  = 19536
    5  = 163
    13 = 62
    14 = 211
    15 = 248
    17 = 1519
    20 = 112
    21 = 248
    22 = 8189
    27 = 1270
    28 = 18
    29 = 14
    31 = 6534
    32 = 132
    33 = 79
    34 = 187
    36 = 308
    37 = 5573
    38 = 1017
    39 = 28
    41 = 62
}

# row -> new F value, for the "全部类型" worksheet
$sheet4Updates = @{
    2  = 15172
    3  = 19536
    5  = 163
    13 = 62
    14 = 211
    15 = 248
    17 = 1519
    21 = 112
    22 = 248
    23 = 8189
    28 = 1270
    29 = 18
    30 = 14
    34 = 6534
    35 = 132
    36 = 79
    37 = 187
    39 = 308
    40 = 5573
    41 = 1017
    42 = 28
    44 = 62
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
